$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D).
# The new column carries over the row styles automatically.
$ws.Columns.Item(2).Insert()

# Row 1 only ever had a value in column A; the column insert stamps an
# (empty) styled cell into B1 too, which the source workbook never had.
# Drop it completely so B1 stays a true empty cell.
$ws.Range("B1").Clear()

# Simple text values (no literal TRUE/FALSE ambiguity) can be set directly.
$ws.Range("B2").Value = "Internal"
$ws.Range("B4").Value = "Internal"

# The word FALSE would normally be auto-converted to a native boolean by
# the smart-typing Value setter (exactly like typing FALSE into a real
# Excel cell). We need it stored as literal text instead, so: write it
# with a harmless trailing space (keeps it text), clean it up with a
# helper formula, then paste back just the computed text value. This
# avoids both the boolean coercion and the "quote prefix" style Excel
# would otherwise stamp onto the cell if we used a leading apostrophe.
function Set-LiteralText($range, [string]$text) {
    $ws.Range($range).Value = ($text + " ")
    $helper = $ws.Range("Z100")
    $helper.Formula = "=TRIM(" + $range + ")"
    $helper.Copy()
    $ws.Range($range).PasteSpecial(-4163)  # xlPasteValues
    $helper.ClearContents()
}

Set-LiteralText "B3" "FALSE"
Set-LiteralText "B5" "FALSE"
Set-LiteralText "B6" "FALSE"

$excel.CutCopyMode = $false

# Match the post-edit selection recorded in the workbook.
$ws.Range("B4:B6").Select() | Out-Null
